$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.714.00"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.913.07"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'239.71"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "'0.4937"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").Value = "'0.2968"
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("D9").Value = "'0.06766"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").Value = "1.923.74"
$ws.Range("E10").Value = "  +1.49%  "
$ws.Range("D11").Value = "'17.08"
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").Value = "'0.07367"
$ws.Range("E12").Value = "  +1.74%  "
$ws.Range("D13").Value = "'5.172"
$ws.Range("E13").Value = "  +2.75%  "
$ws.Range("D14").Value = "'88.37"
$ws.Range("E14").Value = "  -2.67%  "
$ws.Range("D15").Value = "'0.6729"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").Value = "30.681.14"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "'0.000007920"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").Value = "'13.51"
$ws.Range("E18").Value = "  +3.16%  "
$ws.Range("D19").Value = "'1.003"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").Value = "2.172.31"
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("D21").Value = "'5.368"
$ws.Range("E21").Value = "  +11.72%  "
$ws.Range("D22").Value = "'1.004"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").Value = "'197.50"
$ws.Range("E23").Value = "  +2.70%  "
$ws.Range("D24").Value = "'6.260"
$ws.Range("E24").Value = "  +2.93%  "
$ws.Range("D25").Value = "'9.673"
$ws.Range("E25").Value = "  +3.15%  "
$ws.Range("D26").Value = "'164.00"
$ws.Range("E26").Value = "  +5.02%  "
$ws.Range("D27").Value = "'18.62"
$ws.Range("E27").Value = "  -2.21%  "
$ws.Range("D28").Value = "'1.946"
$ws.Range("E28").Value = "  +2.49%  "
$ws.Range("D29").Value = "'1.495"
$ws.Range("E29").Value = "  +6.54%  "
$ws.Range("D30").Value = "'4.354"
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("D31").Value = "'0.09137"
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("D32").Value = "'4.040"
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("D33").Value = "'0.05231"
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("D34").Value = "'0.7405"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").Value = "'1.111"
$ws.Range("E35").Value = "  +0.40%  "
$ws.Range("D36").Value = "'2.727"
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("D37").Value = "'0.01822"
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("D38").Value = "'2.717"
$ws.Range("E38").Value = "  +1.43%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'0.9222"
$ws.Range("E39").Value = "  -0.84%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'2.076"
$ws.Range("E40").Value = "  -1.96%  "
$ws.Range("D41").Value = "'75.30"
$ws.Range("E41").Value = "  +30.64%  "
$ws.Range("D42").Value = "'0.4449"
$ws.Range("E42").Value = "  +1.29%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'106.86"
$ws.Range("E43").Value = "  +1.63%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.920"
$ws.Range("E44").Value = "  +3.20%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").Value = "'0.1388"
$ws.Range("E46").Value = "  +2.81%  "
$ws.Range("D47").Value = "'7.591"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "'35.50"
$ws.Range("E48").Value = "  +5.62%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.045"
$ws.Range("E49").Value = "  +3.29%  "
$ws.Range("D50").Value = "'0.05875"
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("D51").Value = "'0.4018"
$ws.Range("E51").Value = "  +2.32%  "
